# Auto-generated edit script applying the crypto price/volume update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to Text format first,
# otherwise Excel auto-converts the assigned string into a Number cell -
# the source data stores every value (price/volume/coin/link) as text.
$textForcedCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D17",
    "D20",
    "D21",
    "D23",
    "D27",
    "D28",
    "D31",
    "D33",
    "D34",
    "D35",
    "D40",
    "D42",
    "D43",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($ref in $textForcedCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.464.04"
$ws.Range("E2").Value = "  +3.91%  "
$ws.Range("D3").Value = "2.067.15"
$ws.Range("E3").Value = "  +6.24%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "235.70"
$ws.Range("E5").Value = "  +3.47%  "
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +4.52%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "57.78"
$ws.Range("E7").Value = "  +10.04%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +5.41%  "
$ws.Range("D10").Value = "57.98"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").Value = "0.0760"
$ws.Range("E11").Value = "  +4.72%  "
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  +4.84%  "
$ws.Range("D13").Value = "2.369.73"
$ws.Range("E13").Value = "  +6.02%  "
$ws.Range("D14").Value = "14.25"
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").Value = "20.83"
$ws.Range("E15").Value = "  +8.33%  "
$ws.Range("D16").Value = "0.775"
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("D17").Value = "5.17"
$ws.Range("E17").Value = "  +4.91%  "
$ws.Range("D18").Value = "2.067.72"
$ws.Range("E18").Value = "  +5.33%  "
$ws.Range("D19").Value = "37.394.47"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  +24.85%  "
$ws.Range("D21").Value = "68.60"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").Value = "0.0₃0811"
$ws.Range("E22").Value = "  +3.90%  "
$ws.Range("D23").Value = "225.21"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("D27").Value = "163.52"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").Value = "8.83"
$ws.Range("E28").Value = "  +6.23%  "
$ws.Range("E29").Value = "  +9.31%  "
$ws.Range("E30").Value = "  +10.90%  "
$ws.Range("D31").Value = "19.17"
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").Value = "2.60"
$ws.Range("E34").Value = "  +16.58%  "
$ws.Range("D35").Value = "0.0626"
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("E36").Value = "  +7.76%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("D40").Value = "5.82"
$ws.Range("E40").Value = "  +17.84%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "4.47"
$ws.Range("E42").Value = "  +33.18%  "
$ws.Range("D43").Value = "0.0956"
$ws.Range("E43").Value = "  +12.87%  "
$ws.Range("D44").Value = "1.466.89"
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("E45").Value = "  +11.78%  "
$ws.Range("E46").Value = "  +7.77%  "
$ws.Range("D47").Value = "16.09"
$ws.Range("E47").Value = "  +11.90%  "
$ws.Range("E48").Value = "  +7.82%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.29"
$ws.Range("E49").Value = "  +9.85%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "1.02"
$ws.Range("E50").Value = "  +6.09%  "
$ws.Range("D51").Value = "2.94"
$ws.Range("E51").Value = "  +4.02%  "
